# Clear the "value" header (E1) and its sample data (E3) from the
# importcostcode sheet, then move the active selection to E6 (mirrors the
# "hide / uncheck" of the value column before the sheet is re-saved).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").ClearContents() | Out-Null
$ws.Range("E3").ClearContents() | Out-Null

$ws.Range("E6").Select() | Out-Null
